# Append: 2025-12-31 18:27 JST
# Two new job postings are inserted at rows 3-4 on the "ランサーズ" sheet,
# the previously-existing row 3 (自動化 job) slides down to row 5, and every
# row's "取得日時" timestamp is refreshed to the new scrape time. Column D
# is also widened to fit the new price text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-12-31 18:27:16"

# Widen the price column (D) to fit the longer values. The stored OOXML
# <col width> is ColumnWidth plus a fixed ~0.8333 padding offset in this
# engine, so back that out to land on exactly 32 in the saved file.
$ws.Columns.Item(4).ColumnWidth = 31.1666666666667

# Remember the row-3 job that is being pushed down to row 5 before it gets
# overwritten by the new rows.
$movedTitle = $ws.Range("B3").Value2
$movedCategory = $ws.Range("C3").Value2
$movedPrice = $ws.Range("D3").Value2
$movedDeadline = $ws.Range("E3").Value2
$movedUrl = $ws.Range("F3").Value2
$movedScore = $ws.Range("G3").Value2
$movedSkills = $ws.Range("H3").Value2

# Refresh the capture timestamp for every existing data row (2 and, after the
# move below, 3 becomes 5 so it is re-stamped along with everything else).
$ws.Range("A2").Value = $timestamp

# Clear existing hyperlinks so they can be rebuilt cleanly in row order.
$ws.Hyperlinks.Delete()

# Row 3: new AI (LLM) engineer job posting.
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5460294"
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# Row 4: new generative-AI (RAG) lead engineer job posting.
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5460267"
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

# Row 5: the job formerly in row 3, moved down, timestamp refreshed.
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = $movedTitle
$ws.Range("C5").Value = $movedCategory
$ws.Range("D5").Value = $movedPrice
$ws.Range("E5").Value = $movedDeadline
$ws.Range("F5").Value = $movedUrl
$ws.Range("G5").Value = $movedScore
$ws.Range("H5").Value = $movedSkills

# Rebuild hyperlinks in row order (F2..F5) so relationship ids come out
# rId1..rId4 matching cell order.
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value2)
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value2)
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value2)
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value2)
$ws.Range("F5").Style = "Hyperlink"
